# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) for rows 2..46, replacing the old Strike# values.
$newK = @{
    2=3; 3=9; 4=6; 5=5; 6=7; 7=3; 8=1; 9=3; 10=4; 11=7; 12=5; 13=2; 14=3; 15=2; 16=1;
    17=2; 18=3; 19=2; 20=0; 21=0; 22=0; 23=3; 24=1; 25=1; 26=2; 27=1; 28=2; 29=0; 30=1;
    31=3; 32=3; 33=2; 34=4; 35=1; 36=2; 37=2; 38=3; 39=3; 40=2; 41=3; 42=6; 43=3; 44=4; 45=5; 46=2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
